# Update "想去人数" (F column) counts on the "展览" and "全部类型" worksheets
# to reflect the latest generated snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 342
$ws1.Range("F4").Value = 36
$ws1.Range("F5").Value = 3402
$ws1.Range("F6").Value = 2148
$ws1.Range("F7").Value = 413
$ws1.Range("F8").Value = 161
$ws1.Range("F9").Value = 45
$ws1.Range("F10").Value = 34
$ws1.Range("F11").Value = 1260
$ws1.Range("F12").Value = 226
$ws1.Range("F13").Value = 1512

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 342
$ws4.Range("F4").Value = 36
$ws4.Range("F5").Value = 3402
$ws4.Range("F6").Value = 2148
$ws4.Range("F7").Value = 413
$ws4.Range("F9").Value = 161
$ws4.Range("F10").Value = 45
$ws4.Range("F11").Value = 34
$ws4.Range("F14").Value = 1260
$ws4.Range("F15").Value = 226
$ws4.Range("F16").Value = 1512

$wb.Save()
